$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2E")
$ws.Rows("31:31").Insert()
$ws.Range("A31").Value2 = "SKU-T28CPB001 -3QTY"
$ws.Range("AI31").Value2 = "28 oz All Around™ Tumbler"
$ws.Range("AJ31").Value2 = "'3"
$ws.Range("AL31").Value2 = "Black"

# Collect hyperlink info first (avoid mutating while iterating)
$count = $ws.Hyperlinks.Count
Write-Host "count=$count"
$toFix = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $count; $i++) {
  $h = $ws.Hyperlinks.Item($i)
  $r = $h.Range
  $row = $r.Row
  $col = $r.Column
  Write-Host ("idx=$i row=$row col=$col")
}
